# Fruta / hortaliza, semanal
# Refresh the weekly Alcachofa (Vega Monumental Concepcion) price rows with
# the latest market report: update existing rows 82-96 and append two new
# observation rows (97-98), pushing the sheet's used range to A1:R98.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 82-96 with refreshed weekly price data ---

# Row 82
$ws.Range("D82").Value = 45127
$ws.Range("H82").Value = 'Argentina(o)'
$ws.Range("J82").Value = 110
$ws.Range("K82").Value = 16000
$ws.Range("L82").Value = 17000
$ws.Range("M82").Value = 16545
$ws.Range("N82").Value = '$/caja 50 unidades'
$ws.Range("P82").Value = 331
$ws.Range("Q82").Value = 50

# Row 83
$ws.Range("D83").Value = 45127
$ws.Range("H83").Value = 'Española'
$ws.Range("J83").Value = 80
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = 20000
$ws.Range("N83").Value = '$/caja 30 unidades'
$ws.Range("P83").Value = 667
$ws.Range("Q83").Value = 30

# Row 84
$ws.Range("D84").Value = 44376
$ws.Range("H84").Value = 'Española'
$ws.Range("K84").Value = 19000
$ws.Range("M84").Value = 19500
$ws.Range("N84").Value = '$/caja 30 unidades'
$ws.Range("P84").Value = 650
$ws.Range("Q84").Value = 30

# Row 85
$ws.Range("D85").Value = 45106
$ws.Range("H85").Value = 'Argentina(o)'
$ws.Range("K85").Value = 14000
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = 14500
$ws.Range("N85").Value = '$/caja 50 unidades'
$ws.Range("P85").Value = 290
$ws.Range("Q85").Value = 50

# Row 86
$ws.Range("D86").Value = 44358
$ws.Range("H86").Value = 'Argentina(o)'
$ws.Range("K86").Value = 18000
$ws.Range("L86").Value = 20000
$ws.Range("M86").Value = 19000
$ws.Range("N86").Value = '$/caja 50 unidades'
$ws.Range("P86").Value = 380
$ws.Range("Q86").Value = 50

# Row 87
$ws.Range("D87").Value = 44358
$ws.Range("H87").Value = 'Española'
$ws.Range("K87").Value = 18000
$ws.Range("M87").Value = 19000
$ws.Range("N87").Value = '$/caja 30 unidades'
$ws.Range("P87").Value = 633
$ws.Range("Q87").Value = 30

# Row 88
$ws.Range("D88").Value = 44420
$ws.Range("K88").Value = 14000
$ws.Range("L88").Value = 15000
$ws.Range("M88").Value = 14500
$ws.Range("P88").Value = 483

# Row 89
$ws.Range("D89").Value = 44364
$ws.Range("J89").Value = 100
$ws.Range("K89").Value = 19000
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = 19500
$ws.Range("N89").Value = '$/caja 50 unidades'
$ws.Range("P89").Value = 390
$ws.Range("Q89").Value = 50

# Row 90
$ws.Range("D90").Value = 44364
$ws.Range("H90").Value = 'Española'
$ws.Range("K90").Value = 19000
$ws.Range("L90").Value = 20000
$ws.Range("M90").Value = 19500
$ws.Range("N90").Value = '$/caja 30 unidades'
$ws.Range("P90").Value = 650
$ws.Range("Q90").Value = 30

# Row 91
$ws.Range("D91").Value = 44729
$ws.Range("J91").Value = 200
$ws.Range("K91").Value = 17000
$ws.Range("L91").Value = 18000
$ws.Range("M91").Value = 17500
$ws.Range("N91").Value = '$/caja 40 unidades'
$ws.Range("P91").Value = 438
$ws.Range("Q91").Value = 40

# Row 92
$ws.Range("D92").Value = 44811
$ws.Range("H92").Value = 'Madrigal'
$ws.Range("J92").Value = 100
$ws.Range("K92").Value = 12000
$ws.Range("L92").Value = 13000
$ws.Range("M92").Value = 12500
$ws.Range("N92").Value = '$/caja 40 unidades'
$ws.Range("P92").Value = 312
$ws.Range("Q92").Value = 40

# Row 93
$ws.Range("D93").Value = 45112
$ws.Range("J93").Value = 100
$ws.Range("K93").Value = 14000
$ws.Range("L93").Value = 15000
$ws.Range("M93").Value = 14500
$ws.Range("P93").Value = 290

# Row 94
$ws.Range("D94").Value = 44749
$ws.Range("J94").Value = 250
$ws.Range("K94").Value = 13000
$ws.Range("L94").Value = 15000
$ws.Range("M94").Value = 14200
$ws.Range("P94").Value = 284

# Row 95
$ws.Range("D95").Value = 45099
$ws.Range("H95").Value = 'Argentina(o)'
$ws.Range("J95").Value = 130
$ws.Range("K95").Value = 16000
$ws.Range("L95").Value = 17000
$ws.Range("M95").Value = 16615
$ws.Range("N95").Value = '$/caja 50 unidades'
$ws.Range("P95").Value = 332
$ws.Range("Q95").Value = 50

# Row 96
$ws.Range("D96").Value = 44741
$ws.Range("H96").Value = 'Argentina(o)'
$ws.Range("J96").Value = 100
$ws.Range("K96").Value = 16000
$ws.Range("L96").Value = 17000
$ws.Range("M96").Value = 16500
$ws.Range("N96").Value = '$/caja 50 unidades'
$ws.Range("P96").Value = 330
$ws.Range("Q96").Value = 50

# --- Append two new rows (97-98) of price data ---

# Row 97
$ws.Range("D97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A97").Value = 11
$ws.Range("B97").Value = 'Vega Monumental Concepción'
$ws.Range("C97").Value = 'Bíobío'
$ws.Range("D97").Value = 44741
$ws.Range("E97").Value = 8
$ws.Range("F97").Value = 100112013
$ws.Range("G97").Value = 'Alcachofa'
$ws.Range("H97").Value = 'Española'
$ws.Range("I97").Value = 'Primera'
$ws.Range("J97").Value = 100
$ws.Range("K97").Value = 20000
$ws.Range("L97").Value = 22000
$ws.Range("M97").Value = 21000
$ws.Range("N97").Value = '$/caja 30 unidades'
$ws.Range("O97").Value = 'Provincia de Limarí'
$ws.Range("P97").Value = 700
$ws.Range("Q97").Value = 30
$ws.Range("R97").Value = 'Hortaliza'

# Row 98
$ws.Range("D98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A98").Value = 11
$ws.Range("B98").Value = 'Vega Monumental Concepción'
$ws.Range("C98").Value = 'Bíobío'
$ws.Range("D98").Value = 44777
$ws.Range("E98").Value = 8
$ws.Range("F98").Value = 100112013
$ws.Range("G98").Value = 'Alcachofa'
$ws.Range("H98").Value = 'Española'
$ws.Range("I98").Value = 'Primera'
$ws.Range("J98").Value = 110
$ws.Range("K98").Value = 18000
$ws.Range("L98").Value = 19000
$ws.Range("M98").Value = 18545
$ws.Range("N98").Value = '$/caja 30 unidades'
$ws.Range("O98").Value = 'Provincia de Limarí'
$ws.Range("P98").Value = 618
$ws.Range("Q98").Value = 30
$ws.Range("R98").Value = 'Hortaliza'

